$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("store")

# Add new headers for latitude / longitude
$ws.Range("E1").Value = "latitude"
$ws.Range("F1").Value = "longitude"

# Fill in latitude/longitude values for each store row
$ws.Range("E2").Value = 3.1390030000000002
$ws.Range("F2").Value = 101.68685499999999

$ws.Range("E3").Value = 3.1341199999999998
$ws.Range("F3").Value = 101.68653

$ws.Range("E4").Value = 3.1352799999999998
$ws.Range("F4").Value = 101.6871

$ws.Range("E5").Value = 3.13428
$ws.Range("F5").Value = 101.68810000000001

$ws.Range("E6").Value = 3.13571
$ws.Range("F6").Value = 101.6961

$ws.Range("E7").Value = 3.1360100000000002
$ws.Range("F7").Value = 101.68899999999999
